$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell B1 from "Date" to "DateTime"
$ws.Range("B1").Value = "DateTime"

# Set a custom width for column B (~18.24 characters)
$ws.Columns.Item(2).ColumnWidth = 17.29

# Add new data row
$ws.Range("A2").Value = "10238201"
$ws.Range("B2").Value = "22/11/2018_09:42:33"

# Move active selection to C4 to match resulting cursor position
$ws.Range("C4").Select()
